$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column AF: "Team" header (copy header formatting from AE1) + "Norwich"
# for every data row (Daniel Farke's club in every one of these fixtures).

$ws.Range("AE1").Copy()
$ws.Range("AF1").PasteSpecial(-4122)
$ws.Range("AF1").Value = "Team"

$ws.Range("AF2:AF131").Value = "Norwich"
